$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> (latitude, longitude) to fill into columns C and D.
$coords = @{
    7  = @(16.4606, -85.82510000000001)
    8  = @(16.44315, -85.8612)
    9  = @(16.4504, -85.8623)
    10 = @(16.4199, -85.90000000000001)
    11 = @(16.5555, -85.89919999999999)
    12 = @(16.4702, -85.90260000000001)
    24 = @(15.909713, -87.930267)
    25 = @(15.887213, -87.936851)
    27 = @(16.393, -86.274)
    34 = @(16.374, -86.283)
    35 = @(16.358, -86.289)
}

foreach ($row in $coords.Keys) {
    $lat = $coords[$row][0]
    $lon = $coords[$row][1]
    $ws.Cells.Item($row, 3).Value = $lat
    $ws.Cells.Item($row, 4).Value = $lon
}
